$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, $val) {
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "62.718.23"
Set-TextValue $ws.Range("E2") "  +2.20%  "
Set-TextValue $ws.Range("D3") "2.945.58"
Set-TextValue $ws.Range("E3") "  +1.00%  "
Set-TextValue $ws.Range("D4") "0.999"
Set-TextValue $ws.Range("E4") "  +0.02%  "
Set-TextValue $ws.Range("D5") "589.37"
Set-TextValue $ws.Range("E5") "  -0.75%  "
Set-TextValue $ws.Range("D6") "147.54"
Set-TextValue $ws.Range("E6") "  +3.07%  "
Set-TextValue $ws.Range("D7") "1.00"
Set-TextValue $ws.Range("E7") "  +0.03%  "
Set-TextValue $ws.Range("D8") "2.949.07"
Set-TextValue $ws.Range("E8") "  +1.12%  "
Set-TextValue $ws.Range("D9") "0.506"
Set-TextValue $ws.Range("E9") "  +1.52%  "
Set-TextValue $ws.Range("D10") "6.94"
Set-TextValue $ws.Range("E10") "  -1.17%  "
Set-TextValue $ws.Range("D11") "0.149"
Set-TextValue $ws.Range("E11") "  +6.10%  "
Set-TextValue $ws.Range("D12") "0.435"
Set-TextValue $ws.Range("E12") "  +0.22%  "
Set-TextValue $ws.Range("D13") "0.0000233"
Set-TextValue $ws.Range("E13") "  +5.33%  "
Set-TextValue $ws.Range("D14") "32.20"
Set-TextValue $ws.Range("E14") "  -1.68%  "
Set-TextValue $ws.Range("D15") "0.126"
Set-TextValue $ws.Range("E15") "  -1.28%  "
Set-TextValue $ws.Range("D16") "3.432.98"
Set-TextValue $ws.Range("E16") "  +1.00%  "
Set-TextValue $ws.Range("D17") "62.688.00"
Set-TextValue $ws.Range("E17") "  +2.26%  "
Set-TextValue $ws.Range("D18") "6.65"
Set-TextValue $ws.Range("E18") "  +0.99%  "
Set-TextValue $ws.Range("D19") "2.948.98"
Set-TextValue $ws.Range("E19") "  +1.12%  "
Set-TextValue $ws.Range("D20") "436.56"
Set-TextValue $ws.Range("E20") "  +1.40%  "
Set-TextValue $ws.Range("D21") "13.40"
Set-TextValue $ws.Range("E21") "  -0.77%  "
Set-TextValue $ws.Range("D22") "0.660"
Set-TextValue $ws.Range("E22") "  -0.75%  "
Set-TextValue $ws.Range("D23") "6.95"
Set-TextValue $ws.Range("E23") "  -0.85%  "
Set-TextValue $ws.Range("D24") "11.23"
Set-TextValue $ws.Range("E24") "  +3.89%  "
Set-TextValue $ws.Range("D25") "80.35"
Set-TextValue $ws.Range("E25") "  -0.87%  "
Set-TextValue $ws.Range("D26") "11.82"
Set-TextValue $ws.Range("E26") "  +1.96%  "
Set-TextValue $ws.Range("D27") "2.10"
Set-TextValue $ws.Range("E27") "  -1.46%  "
Set-TextValue $ws.Range("D28") "1.00"
Set-TextValue $ws.Range("E28") "  +0.00%  "
Set-TextValue $ws.Range("D29") "7.24"
Set-TextValue $ws.Range("E29") "  +5.88%  "
Set-TextValue $ws.Range("D30") "2.19"
Set-TextValue $ws.Range("E30") "  +1.73%  "
Set-TextValue $ws.Range("D31") "2.59"
Set-TextValue $ws.Range("E31") "  +0.58%  "
Set-TextValue $ws.Range("D32") "0.0000100"
Set-TextValue $ws.Range("E32") "  +14.77%  "
Set-TextValue $ws.Range("D33") "0.108"
Set-TextValue $ws.Range("E33") "  +0.98%  "
Set-TextValue $ws.Range("D34") "26.21"
Set-TextValue $ws.Range("E34") "  -1.02%  "
Set-TextValue $ws.Range("D35") "1.00"
Set-TextValue $ws.Range("E35") "  +0.08%  "
Set-TextValue $ws.Range("D36") "0.992"
Set-TextValue $ws.Range("E36") "  -0.69%  "
Set-TextValue $ws.Range("D37") "5.56"
Set-TextValue $ws.Range("E37") "  +0.43%  "
Set-TextValue $ws.Range("D38") "3.02"
Set-TextValue $ws.Range("E38") "  +3.75%  "
Set-TextValue $ws.Range("D39") "49.69"
Set-TextValue $ws.Range("E39") "  +0.39%  "
Set-TextValue $ws.Range("D40") "1.99"
Set-TextValue $ws.Range("E40") "  +2.02%  "
Set-TextValue $ws.Range("D41") "8.36"
Set-TextValue $ws.Range("E41") "  -0.62%  "
Set-TextValue $ws.Range("D42") "0.115"
Set-TextValue $ws.Range("E42") "  -5.23%  "
Set-TextValue $ws.Range("D43") "0.275"
Set-TextValue $ws.Range("E43") "  +0.94%  "
Set-TextValue $ws.Range("D44") "39.82"
Set-TextValue $ws.Range("E44") "  -3.51%  "
Set-TextValue $ws.Range("D45") "134.52"
Set-TextValue $ws.Range("E45") "  +0.86%  "
Set-TextValue $ws.Range("D46") "2.677.71"
Set-TextValue $ws.Range("E46") "  -0.06%  "
Set-TextValue $ws.Range("D47") "0.0335"
Set-TextValue $ws.Range("E47") "  -0.95%  "
Set-TextValue $ws.Range("D48") "352.58"
Set-TextValue $ws.Range("E48") "  +0.98%  "
Set-TextValue $ws.Range("D49") "1.00"
Set-TextValue $ws.Range("E49") "  +0.00%  "
Set-TextValue $ws.Range("D50") "0.104"
Set-TextValue $ws.Range("E50") "  -0.25%  "
Set-TextValue $ws.Range("D51") "22.49"
Set-TextValue $ws.Range("E51") "  -3.06%  "
